$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table held years 2000-2020 in rows 2-22. The update drops the
# 2000-2009 rows (old rows 2-11), keeps 2010-2019 (old rows 12-21,
# becoming new rows 2-11), revises the 2020 figures (new row 12) and
# appends a new 2021 row (new row 13).

# Remove the obsolete 2000-2009 rows; rows below shift up automatically.
$ws.Rows("2:11").Delete() | Out-Null

# Revise the 2020 row (now row 12) with the corrected figures.
$ws.Cells.Item(12, 2).Value = 1.68
$ws.Cells.Item(12, 3).Value = 1.68
$ws.Cells.Item(12, 4).Value = 1
$ws.Cells.Item(12, 5).Value = 1.14

# Add the new 2021 row (row 13), copying row 12's formatting first so the
# year label picks up the same style used by the rest of column A.
$ws.Range("A12").Copy($ws.Range("A13")) | Out-Null
$ws.Cells.Item(13, 1).Value = "2021年"
$ws.Cells.Item(13, 2).Value = 1.28
$ws.Cells.Item(13, 3).Value = 1.2
$ws.Cells.Item(13, 4).Value = 0.64
$ws.Cells.Item(13, 5).Value = 0.77
